# Apply the Apr 19 2023 04:29:24 UTC cryptos-list GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (matches the source data, which
# stores prices/percentages as text, not numbers) without leaving a
# lingering custom number-format on the cell.
function Set-TextValue {
    param($Sheet, $CellRef, $NewValue)
    $cell = $Sheet.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '30.276.26'
Set-TextValue $ws 'E2' '  +1.96%  '
Set-TextValue $ws 'D3' '2.092.67'
Set-TextValue $ws 'E3' '  -0.24%  '
Set-TextValue $ws 'E4' '  -0.59%  '
Set-TextValue $ws 'D5' '341.92'
Set-TextValue $ws 'E5' '  -0.38%  '
Set-TextValue $ws 'D6' '1.002'
Set-TextValue $ws 'E6' '  -0.53%  '
Set-TextValue $ws 'D7' '0.5302'
Set-TextValue $ws 'E7' '  +2.32%  '
Set-TextValue $ws 'D8' '0.4379'
Set-TextValue $ws 'E8' '  +0.18%  '
Set-TextValue $ws 'D9' '54.64'
Set-TextValue $ws 'E9' '  +2.15%  '
Set-TextValue $ws 'D10' '0.09344'
Set-TextValue $ws 'E10' '  +1.39%  '
Set-TextValue $ws 'E11' '  +0.59%  '
Set-TextValue $ws 'D12' '24.65'
Set-TextValue $ws 'E12' '  +0.18%  '
Set-TextValue $ws 'D13' '8.548'
Set-TextValue $ws 'E13' '  +4.77%  '
Set-TextValue $ws 'D14' '6.872'
Set-TextValue $ws 'D15' '2.049.90'
Set-TextValue $ws 'E15' '  -1.53%  '
Set-TextValue $ws 'D16' '101.17'
Set-TextValue $ws 'E16' '  -1.85%  '
Set-TextValue $ws 'E17' '  +0.13%  '
Set-TextValue $ws 'E18' '  -0.54%  '
Set-TextValue $ws 'D19' '21.07'
Set-TextValue $ws 'E19' '  +0.28%  '
Set-TextValue $ws 'D20' '0.06724'
Set-TextValue $ws 'E20' '  +0.80%  '
Set-TextValue $ws 'D21' '6.334'
Set-TextValue $ws 'E21' '  +2.05%  '
Set-TextValue $ws 'D22' '1.001'
Set-TextValue $ws 'E22' '  -0.62%  '
Set-TextValue $ws 'D23' '30.260.87'
Set-TextValue $ws 'E23' '  +1.81%  '
Set-TextValue $ws 'D24' '12.43'
Set-TextValue $ws 'E24' '  -0.97%  '
Set-TextValue $ws 'D25' '2.319'
Set-TextValue $ws 'E25' '  +0.66%  '
Set-TextValue $ws 'D26' '6.952'
Set-TextValue $ws 'E26' '  +9.33%  '
Set-TextValue $ws 'D27' '21.78'
Set-TextValue $ws 'E27' '  -0.59%  '
Set-TextValue $ws 'D28' '162.09'
Set-TextValue $ws 'E28' '  +0.17%  '
Set-TextValue $ws 'D29' '2.500'
Set-TextValue $ws 'E29' '  +0.41%  '
Set-TextValue $ws 'D30' '133.75'
Set-TextValue $ws 'E30' '  +0.15%  '
Set-TextValue $ws 'D31' '1.129'
Set-TextValue $ws 'E31' '  +0.17%  '
Set-TextValue $ws 'D32' '0.1053'
Set-TextValue $ws 'E32' '  +0.20%  '
Set-TextValue $ws 'D33' '1.662'
Set-TextValue $ws 'E33' '  -1.39%  '
Set-TextValue $ws 'D34' '6.236'
Set-TextValue $ws 'E34' '  +0.67%  '
Set-TextValue $ws 'D35' '3.913'
Set-TextValue $ws 'E35' '  -0.96%  '
Set-TextValue $ws 'D36' '10.05'
Set-TextValue $ws 'E36' '  -3.46%  '
Set-TextValue $ws 'D37' '0.02611'
Set-TextValue $ws 'E37' '  +1.48%  '
Set-TextValue $ws 'D38' '0.06744'
Set-TextValue $ws 'E38' '  +0.40%  '
Set-TextValue $ws 'D39' '12.55'
Set-TextValue $ws 'E39' '  +0.52%  '
Set-TextValue $ws 'D40' '0.6947'
Set-TextValue $ws 'E40' '  -0.61%  '
Set-TextValue $ws 'D41' '1.340'
Set-TextValue $ws 'E41' '  +1.17%  '
Set-TextValue $ws 'D42' '0.2207'
Set-TextValue $ws 'E42' '  -0.43%  '
Set-TextValue $ws 'D43' '0.6769'
Set-TextValue $ws 'E43' '  -0.13%  '
Set-TextValue $ws 'B44' 'NEARProtocol'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D44' '2.346'
Set-TextValue $ws 'E44' '  +1.05%  '
Set-TextValue $ws 'B45' 'EnergySwap'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D45' '14.18'
Set-TextValue $ws 'E45' '  -0.51%  '
Set-TextValue $ws 'E46' '  -0.47%  '
Set-TextValue $ws 'E47' '  +7.98%  '
Set-TextValue $ws 'D48' '3.635'
Set-TextValue $ws 'E48' '  +0.33%  '
Set-TextValue $ws 'D49' '0.00000000345'
Set-TextValue $ws 'E49' '  -3.87%  '
Set-TextValue $ws 'D50' '1.208'
Set-TextValue $ws 'E50' '  +4.68%  '
Set-TextValue $ws 'D51' '1.212'
Set-TextValue $ws 'E51' '  -0.38%  '

Write-Output "Applied 97 cell updates."
